# Update the 380 kV case results: loading_percent values for rows 2-25
# (columns B, C, E, F, G, H, I, J, O). Columns D, K, L, M, N and column A
# (the index) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.04611075275231
$ws.Range("C2").Value = 10.96054116066838
$ws.Range("E2").Value = 26.3633038781184
$ws.Range("F2").Value = 37.25174455987145
$ws.Range("G2").Value = 16.46401062349421
$ws.Range("H2").Value = 11.31324375025298
$ws.Range("I2").Value = 15.21513753360844
$ws.Range("J2").Value = 7.140670539969022
$ws.Range("O2").Value = 15.55679373266467
$ws.Range("B3").Value = 14.18553567165188
$ws.Range("C3").Value = 10.35457268951409
$ws.Range("E3").Value = 26.03240471993254
$ws.Range("F3").Value = 37.07414304041676
$ws.Range("G3").Value = 16.65597903237424
$ws.Range("H3").Value = 11.38593106092411
$ws.Range("I3").Value = 15.37345971912739
$ws.Range("J3").Value = 7.167891829517322
$ws.Range("O3").Value = 15.69172910135928
$ws.Range("B4").Value = 13.62806503764386
$ws.Range("C4").Value = 9.962567813385331
$ws.Range("E4").Value = 25.83240447746792
$ws.Range("F4").Value = 36.97670805082707
$ws.Range("G4").Value = 16.78543430866813
$ws.Range("H4").Value = 11.43314453608169
$ws.Range("I4").Value = 15.47543336176097
$ws.Range("J4").Value = 7.185518762329643
$ws.Range("O4").Value = 15.77997731378135
$ws.Range("B5").Value = 13.39370027933022
$ws.Range("C5").Value = 9.797906278108581
$ws.Range("E5").Value = 25.75178950003703
$ws.Range("F5").Value = 36.93995253649835
$ws.Range("G5").Value = 16.84106179735047
$ws.Range("H5").Value = 11.45303417820285
$ws.Range("I5").Value = 15.51818977108192
$ws.Range("J5").Value = 7.192931963521472
$ws.Range("O5").Value = 15.81729239452477
$ws.Range("B6").Value = 13.35435393859508
$ws.Range("C6").Value = 9.770270681064556
$ws.Range("E6").Value = 25.73845956494397
$ws.Range("F6").Value = 36.93402826181288
$ws.Range("G6").Value = 16.85047109001479
$ws.Range("H6").Value = 11.45637608938344
$ws.Range("I6").Value = 15.52536210256142
$ws.Range("J6").Value = 7.194176831557439
$ws.Range("O6").Value = 15.82357015306349
$ws.Range("B7").Value = 13.62493324942157
$ws.Range("C7").Value = 9.960366892571901
$ws.Range("E7").Value = 25.83131356919707
$ws.Range("F7").Value = 36.97620037284013
$ws.Range("G7").Value = 16.78617293851086
$ws.Range("H7").Value = 11.43341014357608
$ws.Range("I7").Value = 15.47600512080691
$ws.Range("J7").Value = 7.185617807073188
$ws.Range("O7").Value = 15.7804750840548
$ws.Range("B8").Value = 14.75548245479289
$ws.Range("C8").Value = 10.75579330810931
$ws.Range("E8").Value = 26.24860499705185
$ws.Range("F8").Value = 37.18811664537455
$ws.Range("G8").Value = 16.52777247644767
$ws.Range("H8").Value = 11.33777031701402
$ws.Range("I8").Value = 15.26874058837971
$ws.Range("J8").Value = 7.149867267390371
$ws.Range("O8").Value = 15.60219678386414
$ws.Range("B9").Value = 16.73779829651475
$ws.Range("C9").Value = 12.15412482338277
$ws.Range("E9").Value = 27.08812788675933
$ws.Range("F9").Value = 37.69428798973411
$ws.Range("G9").Value = 16.11479052181018
$ws.Range("H9").Value = 11.17070775388304
$ws.Range("I9").Value = 14.89992395063834
$ws.Range("J9").Value = 7.086980872810811
$ws.Range("O9").Value = 15.29560352597122
$ws.Range("B10").Value = 18.04720459732215
$ws.Range("C10").Value = 13.07970375981681
$ws.Range("E10").Value = 27.71251708452303
$ws.Range("F10").Value = 38.11902850127653
$ws.Range("G10").Value = 15.87102835737547
$ws.Range("H10").Value = 11.06044319924986
$ws.Range("I10").Value = 14.65165990312714
$ws.Range("J10").Value = 7.045146842234215
$ws.Range("O10").Value = 15.09683182604299
$ws.Range("B11").Value = 18.61045508113943
$ws.Range("C11").Value = 13.47822004851101
$ws.Range("E11").Value = 27.997088419508
$ws.Range("F11").Value = 38.32313520087832
$ws.Range("G11").Value = 15.77363317837518
$ws.Range("H11").Value = 11.01298831120167
$ws.Range("I11").Value = 14.54359864880646
$ws.Range("J11").Value = 7.02705708993967
$ws.Range("O11").Value = 15.01221681444934
$ws.Range("B12").Value = 18.81905531270576
$ws.Range("C12").Value = 13.62586117135335
$ws.Range("E12").Value = 28.10483215457156
$ws.Range("F12").Value = 38.40193516945716
$ws.Range("G12").Value = 15.73873765185682
$ws.Range("H12").Value = 10.99540742441068
$ws.Range("I12").Value = 14.50337602320876
$ws.Range("J12").Value = 7.020341743093763
$ws.Range("O12").Value = 14.98101597603151
$ws.Range("B13").Value = 18.77433852308106
$ws.Range("C13").Value = 13.59420974153138
$ws.Range("E13").Value = 28.08162985553969
$ws.Range("F13").Value = 38.3848980042543
$ws.Range("G13").Value = 15.7461639792655
$ws.Range("H13").Value = 10.9991764660111
$ws.Range("I13").Value = 14.51200770050189
$ws.Range("J13").Value = 7.021782022474383
$ws.Range("O13").Value = 14.98769813213881
$ws.Range("B14").Value = 18.62771089543075
$ws.Range("C14").Value = 13.49043220895106
$ws.Range("E14").Value = 28.00595343835118
$ws.Range("F14").Value = 38.32958821244942
$ws.Range("G14").Value = 15.77072224104666
$ws.Range("H14").Value = 11.01153411816053
$ws.Range("I14").Value = 14.54027554214458
$ws.Range("J14").Value = 7.026501914158183
$ws.Range("O14").Value = 15.00963301214495
$ws.Range("B15").Value = 18.53728564541795
$ws.Range("C15").Value = 13.42643913582009
$ws.Range("E15").Value = 27.95959444571434
$ws.Range("F15").Value = 38.29590416635626
$ws.Range("G15").Value = 15.78602486113423
$ws.Range("H15").Value = 11.01915423775706
$ws.Range("I15").Value = 14.5576811934044
$ws.Range("J15").Value = 7.029410532888922
$ws.Range("O15").Value = 15.02317846772113
$ws.Range("B16").Value = 18.00973993993271
$ws.Range("C16").Value = 13.05320377376299
$ws.Range("E16").Value = 27.69392233753337
$ws.Range("F16").Value = 38.10590400832146
$ws.Range("G16").Value = 15.87766878109417
$ws.Range("H16").Value = 11.06359894145612
$ws.Range("I16").Value = 14.65881978643286
$ws.Range("J16").Value = 7.046347943596007
$ws.Range("O16").Value = 15.10247897568053
$ws.Range("B17").Value = 17.67778310590272
$ws.Range("C17").Value = 12.81844206250405
$ws.Range("E17").Value = 27.53101030210025
$ws.Range("F17").Value = 37.99209497308571
$ws.Range("G17").Value = 15.93737779046061
$ws.Range("H17").Value = 11.09155732107206
$ws.Range("I17").Value = 14.7221114071281
$ws.Range("J17").Value = 7.056979138591172
$ws.Range("O17").Value = 15.15261879639971
$ws.Range("B18").Value = 17.48380136130296
$ws.Range("C18").Value = 12.68129344220677
$ws.Range("E18").Value = 27.43736334664263
$ws.Range("F18").Value = 37.92766368388038
$ws.Range("G18").Value = 15.97298817162604
$ws.Range("H18").Value = 11.10789283005466
$ws.Range("I18").Value = 14.7589741516594
$ws.Range("J18").Value = 7.063182494566678
$ws.Range("O18").Value = 15.18200412679897
$ws.Range("B19").Value = 17.41759990802797
$ws.Range("C19").Value = 12.63449415983349
$ws.Range("E19").Value = 27.40566856434466
$ws.Range("F19").Value = 37.90602670366069
$ws.Range("G19").Value = 15.98526151214902
$ws.Range("H19").Value = 11.11346747121408
$ws.Range("I19").Value = 14.77153418989463
$ws.Range("J19").Value = 7.065298070188154
$ws.Range("O19").Value = 15.19204713336427
$ws.Range("B20").Value = 17.71343628494787
$ws.Range("C20").Value = 12.84365247472219
$ws.Range("E20").Value = 27.54834748457822
$ws.Range("F20").Value = 38.00410407750842
$ws.Range("G20").Value = 15.93089016937915
$ws.Range("H20").Value = 11.08855475307515
$ws.Range("I20").Value = 14.71532641669293
$ws.Range("J20").Value = 7.055838265817223
$ws.Range("O20").Value = 15.14722475160665
$ws.Range("B21").Value = 18.67090648174124
$ws.Range("C21").Value = 13.52100307401846
$ws.Range("E21").Value = 28.02818264371713
$ws.Range("F21").Value = 38.3457935491705
$ws.Range("G21").Value = 15.76345462441722
$ws.Range("H21").Value = 11.00789380959336
$ws.Range("I21").Value = 14.53195367866847
$ws.Range("J21").Value = 7.025111911163284
$ws.Range("O21").Value = 15.00316733208366
$ws.Range("B22").Value = 19.26932069822082
$ws.Range("C22").Value = 13.9446324007388
$ws.Range("E22").Value = 28.34164002044952
$ws.Range("F22").Value = 38.57787547860512
$ws.Range("G22").Value = 15.6656222713431
$ws.Range("H22").Value = 10.95744623181193
$ws.Range("I22").Value = 14.41617497939826
$ws.Range("J22").Value = 7.005816275852853
$ws.Range("O22").Value = 14.91392259948948
$ws.Range("B23").Value = 18.95244593412907
$ws.Range("C23").Value = 13.72028468892436
$ws.Range("E23").Value = 28.17438547521919
$ws.Range("F23").Value = 38.45322630834846
$ws.Range("G23").Value = 15.71676114245119
$ws.Range("H23").Value = 10.98416333362951
$ws.Range("I23").Value = 14.47759727352392
$ws.Range("J23").Value = 7.016042957577954
$ws.Range("O23").Value = 14.96110333764438
$ws.Range("B24").Value = 17.6973272461799
$ws.Range("C24").Value = 12.83226163850307
$ws.Range("E24").Value = 27.54050929558542
$ws.Range("F24").Value = 37.99867164508593
$ws.Range("G24").Value = 15.93381923021146
$ws.Range("H24").Value = 11.08991139836833
$ws.Range("I24").Value = 14.71839242927098
$ws.Range("J24").Value = 7.056353769841622
$ws.Range("O24").Value = 15.14966165753658
$ws.Range("B25").Value = 16.22712629304014
$ws.Range("C25").Value = 11.79352313188887
$ws.Range("E25").Value = 26.85929520980517
$ws.Range("F25").Value = 37.54788349597856
$ws.Range("G25").Value = 16.21622010880894
$ws.Range("H25").Value = 11.21371088856975
$ws.Range("I25").Value = 14.99569451517753
$ws.Range("J25").Value = 7.103223774195043
$ws.Range("O25").Value = 15.37391439266611
